# Reg_data.xlsx rework: expand the registration test-data sheet from the
# old Name/Password/Email/Cloth-Category layout to a wider
# Name/Email/Password/Company/Phone/Address/Cloth-Category layout, add a
# hyperlink on the email cell, new column widths, and switch the workbook
# to manual calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- calculation mode --------------------------------------------------------
$excel.Calculation = -4135   # xlCalculationManual

# --- drop the legacy "Text" number-format styling that used to live on
#     columns G:H (numFmtId 49 + left/center alignment) before we repurpose
#     those columns -------------------------------------------------------
$ws.Columns.Item(7).ClearFormats()
$ws.Columns.Item(8).ClearFormats()

# --- clear the old data, and drop the old 3rd row entirely (the new layout
#     only needs a header row + one data row) --------------------------------
$ws.Cells.ClearContents()
$ws.Rows.Item(3).Delete()

# --- headers (row 1) ----------------------------------------------------------
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Last Name"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Confirm Password"
$ws.Range("F1").Value = "Company"
$ws.Range("G1").Value = "Phone Number"
$ws.Range("H1").Value = "Street Address"
$ws.Range("I1").Value = "City"
$ws.Range("J1").Value = "State"
$ws.Range("K1").Value = "Zip Code"
$ws.Range("L1").Value = "Cloth Category"
$ws.Range("M1").Value = "Size"
$ws.Range("N1").Value = "Color"
$ws.Range("A1:N1").Font.Bold = $true

# --- data (row 2) --------------------------------------------------------------
$ws.Range("A2").Value = "Rafael"
$ws.Range("B2").Value = "Buckridge"
$ws.Range("C2").Value = "Nazz72@yopmail.com"
$ws.Range("D2").Value = "szMYqE9k@1"
$ws.Range("E2").Value = "szMYqE9k@1"
$ws.Range("F2").Value = "Nazneen's Software Company"
$ws.Range("G2").Value = 1546745564
$ws.Range("H2").Value = "Gulshan, Police Plaza"
$ws.Range("I2").Value = "Dhaka"
$ws.Range("J2").Value = "Alaska"
$ws.Range("K2").Value = 1215
$ws.Range("L2").Value = "Women_Bottoms_Pants"
$ws.Range("M2").Value = 28
$ws.Range("N2").Value = "Blue"

# --- hyperlink on the email cell (keep the visible text as the email, the
#     hyperlink's stored display text is the distinct "ahsGmLIT@1" string) ---
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ahsGmLIT@1", [System.Type]::Missing, [System.Type]::Missing, "ahsGmLIT@1") | Out-Null
$ws.Range("C2").Value = "Nazz72@yopmail.com"
$ws.Range("C2").Style = "Normal"

# --- column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.608072916666666
$ws.Columns.Item(2).ColumnWidth = 11.608072916666666
$ws.Columns.Item(3).ColumnWidth = 20.385416666666668
$ws.Columns.Item(4).ColumnWidth = 27.053385416666668
$ws.Columns.Item(5).ColumnWidth = 17.608072916666668
$ws.Columns.Item(6).ColumnWidth = 27.721354166666668
$ws.Columns.Item(7).ColumnWidth = 23.830729166666668
$ws.Columns.Item(8).ColumnWidth = 19.053385416666668
$ws.Columns.Item(11).ColumnWidth = 14.385416666666666
$ws.Columns.Item(12).ColumnWidth = 21.166666666666668
$ws.Columns.Item(13).ColumnWidth = 16.608072916666668

# --- selection ---------------------------------------------------------------
$ws.Range("D11").Select()
